$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(54).Insert()

$ws.Range("A54").Value = 2
$ws.Range("B54").Value = 'Comercializadora del Agro de Limarí'
$ws.Range("C54").Value = 'Coquimbo'
$ws.Range("D54").Value = 44413
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = 100112021
$ws.Range("G54").Value = 'Ají'
$ws.Range("H54").Value = 'Americana (o)'
$ws.Range("I54").Value = 'Primera'
$ws.Range("J54").Value = 160
$ws.Range("K54").Value = 65000
$ws.Range("L54").Value = 70000
$ws.Range("M54").Value = 67500
$ws.Range("N54").Value = '$/caja 25 kilos'
$ws.Range("O54").Value = 'Provincia de Limarí'
$ws.Range("P54").Value = 2700
$ws.Range("Q54").Value = 25
$ws.Range("R54").Value = 'Hortaliza'
